$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.971.88"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "3.097.32"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'522.88"
$ws.Range("D6").Value = "'144.16"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").Value = "'7.39"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'0.110"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("D12").Value = "3.628.75"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "'26.84"
$ws.Range("E14").Value = "  +4.72%  "
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").Value = "58.941.50"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").Value = "3.093.28"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "'6.16"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'12.97"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'8.13"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").Value = "'344.65"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "'65.72"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "'0.172"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").Value = "0.0₃0926"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'6.67"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("D29").Value = "'7.26"
$ws.Range("E29").Value = "  +3.23%  "
$ws.Range("D30").Value = "'1.85"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "'1.21"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").Value = "'21.04"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "'4.65"
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("D35").Value = "'6.15"
$ws.Range("E35").Value = "  +4.54%  "
$ws.Range("D36").Value = "'26.84"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "'1.30"
$ws.Range("E37").Value = "  +5.32%  "
$ws.Range("D38").Value = "'0.0687"
$ws.Range("D39").Value = "'3.94"
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").Value = "3.136.33"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "'36.77"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'0.666"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  +6.12%  "
$ws.Range("D45").Value = "2.286.27"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'0.0256"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "'20.93"
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("D48").Value = "'0.967"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("D50").Value = "'0.757"
$ws.Range("E50").Value = "  +9.91%  "
$ws.Range("D51").Value = "'261.63"
$ws.Range("E51").Value = "  +10.50%  "
